$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2334293948126801
$ws.Range("C2").Value = 0.4783861671469741
$ws.Range("J2").Value = 0.008645533141210375
$ws.Range("P2").Value = 0.1585014409221902
$ws.Range("S2").Value = 0.1210374639769452
$ws.Range("B3").Value = 0.03333333333333333
$ws.Range("C3").Value = 0.04444444444444445
$ws.Range("J3").Value = 0.03333333333333333
$ws.Range("P3").Value = 0.7055555555555556
$ws.Range("S3").Value = 0.1833333333333333
$ws.Range("P4").Value = 0.78125
$ws.Range("S4").Value = 0.21875
$ws.Range("B6").Value = 0.06024096385542169
$ws.Range("D6").Value = 0.01606425702811245
$ws.Range("F6").Value = 0.08433734939759036
$ws.Range("J6").Value = 0.2650602409638554
$ws.Range("O6").Value = 0.01204819277108434
$ws.Range("Q6").Value = 0.0963855421686747
$ws.Range("R6").Value = 0.07228915662650602
$ws.Range("S6").Value = 0.393574297188755
$ws.Range("B7").Value = 0.1171171171171171
$ws.Range("D7").Value = 0.01801801801801802
$ws.Range("F7").Value = 0.07657657657657657
$ws.Range("J7").Value = 0.1306306306306306
$ws.Range("O7").Value = 0.01351351351351351
$ws.Range("Q7").Value = 0.1441441441441441
$ws.Range("R7").Value = 0.08108108108108109
$ws.Range("S7").Value = 0.4189189189189189
$ws.Range("B8").Value = 0.07954545454545454
$ws.Range("D8").Value = 0.01363636363636364
$ws.Range("E8").Value = 0.004545454545454545
$ws.Range("F8").Value = 0.04772727272727273
$ws.Range("J8").Value = 0.15
$ws.Range("O8").Value = 0.02045454545454545
$ws.Range("Q8").Value = 0.2
$ws.Range("R8").Value = 0.1
$ws.Range("S8").Value = 0.3840909090909091
$ws.Range("B9").Value = 0.0995260663507109
$ws.Range("D9").Value = 0.009478672985781991
$ws.Range("F9").Value = 0.05687203791469194
$ws.Range("J9").Value = 0.1374407582938389
$ws.Range("O9").Value = 0.03317535545023697
$ws.Range("Q9").Value = 0.1943127962085308
$ws.Range("R9").Value = 0.04739336492890995
$ws.Range("S9").Value = 0.4218009478672986
$ws.Range("B10").Value = 0.1119718309859155
$ws.Range("D10").Value = 0.01267605633802817
$ws.Range("E10").Value = 0.0007042253521126761
$ws.Range("F10").Value = 0.06901408450704226
$ws.Range("J10").Value = 0.1387323943661972
$ws.Range("O10").Value = 0.01197183098591549
$ws.Range("Q10").Value = 0.2014084507042254
$ws.Range("R10").Value = 0.07464788732394366
$ws.Range("S10").Value = 0.3788732394366197
$ws.Range("G11").Value = 0.1368715083798883
$ws.Range("J11").Value = 0.08659217877094973
$ws.Range("K11").Value = 0.2067039106145251
$ws.Range("L11").Value = 0.5474860335195531
$ws.Range("S11").Value = 0.0223463687150838
$ws.Range("G12").Value = 0.7102803738317757
$ws.Range("J12").Value = 0.1495327102803738
$ws.Range("K12").Value = 0.01401869158878505
$ws.Range("L12").Value = 0.06542056074766354
$ws.Range("S12").Value = 0.06074766355140187
$ws.Range("G13").Value = 0.717948717948718
$ws.Range("J13").Value = 0.2564102564102564
$ws.Range("S13").Value = 0.02564102564102564
$ws.Range("F15").Value = 0.04508196721311476
$ws.Range("H15").Value = 0.1434426229508197
$ws.Range("I15").Value = 0.06967213114754098
$ws.Range("J15").Value = 0.3811475409836065
$ws.Range("K15").Value = 0.09016393442622951
$ws.Range("M15").Value = 0.01639344262295082
$ws.Range("N15").Value = 0.004098360655737705
$ws.Range("O15").Value = 0.06557377049180328
$ws.Range("S15").Value = 0.1844262295081967
$ws.Range("F16").Value = 0.0196078431372549
$ws.Range("H16").Value = 0.1715686274509804
$ws.Range("I16").Value = 0.05392156862745098
$ws.Range("J16").Value = 0.4558823529411765
$ws.Range("K16").Value = 0.09803921568627451
$ws.Range("M16").Value = 0.0196078431372549
$ws.Range("O16").Value = 0.04901960784313725
$ws.Range("S16").Value = 0.1323529411764706
$ws.Range("F17").Value = 0.02150537634408602
$ws.Range("H17").Value = 0.1763440860215054
$ws.Range("I17").Value = 0.07956989247311828
$ws.Range("J17").Value = 0.4
$ws.Range("K17").Value = 0.1161290322580645
$ws.Range("M17").Value = 0.01720430107526882
$ws.Range("N17").Value = 0.002150537634408602
$ws.Range("O17").Value = 0.07741935483870968
$ws.Range("S17").Value = 0.1096774193548387
$ws.Range("F18").Value = 0.01030927835051546
$ws.Range("H18").Value = 0.1185567010309278
$ws.Range("I18").Value = 0.1185567010309278
$ws.Range("J18").Value = 0.4072164948453608
$ws.Range("K18").Value = 0.1134020618556701
$ws.Range("M18").Value = 0.02061855670103093
$ws.Range("O18").Value = 0.1134020618556701
$ws.Range("S18").Value = 0.09793814432989691
$ws.Range("F19").Value = 0.02295552367288379
$ws.Range("H19").Value = 0.1915351506456241
$ws.Range("I19").Value = 0.08895265423242468
$ws.Range("J19").Value = 0.3715925394548063
$ws.Range("K19").Value = 0.1133428981348637
$ws.Range("M19").Value = 0.01506456241032999
$ws.Range("N19").Value = 0.002869440459110474
$ws.Range("O19").Value = 0.06527977044476327
$ws.Range("S19").Value = 0.1284074605451937
